$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ftests")

# Copy the formatting of the row above (row 40) down into the new rows
# first, so the values we set below land with matching styles.
$ws.Range("B40:I40").Copy()
$ws.Range("B41:I41").PasteSpecial(-4122)
$ws.Range("B40:I40").Copy()
$ws.Range("B42:I42").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 41: fill in the previously blank row
$ws.Range("B41").Value = "fm36"
$ws.Range("C41").Value = "Reverse Franchise deductible calcrule 20"
$ws.Range("D41").Value = 0
$ws.Range("E41").Value = 20
$ws.Range("F41").Value = 1
$ws.Range("G41").Value = 1
$ws.Range("H41").Value = "complete"
$ws.Range("I41").Value = "complete"

# Row 42: brand new row
$ws.Range("B42").Value = "fm37"
$ws.Range("C42").Value = "WE11 Quota Share with % placed and % treaty share for 2 Reinsurers"
$ws.Range("D42").Value = "0,2"
$ws.Range("E42").Value = 22
$ws.Range("F42").Value = 2
$ws.Range("G42").Value = 2
$ws.Range("H42").Value = "complete"
$ws.Range("I42").Value = "complete"

# Match the saved view: scrolled so row 24 is visible, new row selected
$ws.Activate()
$ws.Range("A24").Select()
$ws.Range("C42").Select()
